$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was 112164561 / 93158 -> becomes 112164607 / 93304, Neckera -> Alleniella complanata
$ws.Range("A2").Value = 112164607
$ws.Range("B2").Value = 93304
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 2667
$ws.Range("F2").Value = "Platt fjädermossa"
$ws.Range("G2").Value = "Alleniella complanata"
$ws.Range("H2").Value = "(Hedw.) S.Olsson, Enroth & D.Quandt"
$ws.Range("Q2").Value = 332973
$ws.Range("R2").Value = 6627007

# Row 3: was 112164609 / 92683 (Mnium stellare) -> becomes 112164661 / 89998 (Phlebia serialis)
$ws.Range("A3").Value = 112164661
$ws.Range("B3").Value = 89998
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5467
$ws.Range("F3").Value = "Kådvaxskinn"
$ws.Range("G3").Value = "Phlebia serialis"
$ws.Range("H3").Value = "(Fr.:Fr.) Donk"
$ws.Range("L3").ClearContents()
$ws.Range("Q3").Value = 332865
$ws.Range("R3").Value = 6626972
$ws.Range("AC3").Value = "På granlåga"

# Row 4: was 112164702 / 89369 (Fuscoporia viticola) -> becomes 112164673 / 93303 (Alleniella besseri)
$ws.Range("A4").Value = 112164673
$ws.Range("B4").Value = 93303
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1078
$ws.Range("F4").Value = "Rundfjädermossa"
$ws.Range("G4").Value = "Alleniella besseri"
$ws.Range("H4").Value = "(Lobarz.) S.Olsson, Enroth & D.Quandt"
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 332854
$ws.Range("R4").Value = 6626968
$ws.Range("AC4").Value = "Under överhängande klippa"

# Row 5: was 112164565 / 92683 (Mnium stellare) -> becomes 112164579 / 93308 (Neckera crispa)
$ws.Range("A5").Value = 112164579
$ws.Range("B5").Value = 93308
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2666
$ws.Range("F5").Value = "Grov fjädermossa"
$ws.Range("G5").Value = "Neckera crispa"
$ws.Range("H5").Value = "Hedw."
$ws.Range("Q5").Value = 332923
$ws.Range("R5").Value = 6626955

# Row 6: was 112164673 / 93157 (Neckera besseri) -> becomes 112164561 / 93304 (Alleniella complanata)
$ws.Range("A6").Value = 112164561
$ws.Range("B6").Value = 93304
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 2667
$ws.Range("F6").Value = "Platt fjädermossa"
$ws.Range("G6").Value = "Alleniella complanata"
$ws.Range("H6").Value = "(Hedw.) S.Olsson, Enroth & D.Quandt"
$ws.Range("Q6").Value = 332935
$ws.Range("R6").Value = 6626957
$ws.Range("AC6").ClearContents()

# Row 7: was 112164607 / 93158 (Neckera complanata) -> becomes 112164702 / 89503 (Fuscoporia viticola)
$ws.Range("A7").Value = 112164702
$ws.Range("B7").Value = 89503
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 5447
$ws.Range("F7").Value = "Vedticka"
$ws.Range("G7").Value = "Fuscoporia viticola"
$ws.Range("H7").Value = "(Schwein.) Murrill"
$ws.Range("L7").ClearContents()
$ws.Range("Q7").Value = 332980
$ws.Range("R7").Value = 6627033

# Row 8: was 112164661 / 89864 (Phlebia serialis) -> becomes 112164609 / 92831 (Mnium stellare)
$ws.Range("A8").Value = 112164609
$ws.Range("B8").Value = 92831
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 2362
$ws.Range("F8").Value = "Blek stjärnmossa"
$ws.Range("G8").Value = "Mnium stellare"
$ws.Range("H8").Value = "Hedw."
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 332973
$ws.Range("R8").Value = 6627007
$ws.Range("AC8").ClearContents()

# Row 9: was 112164579 / 93159 (Neckera crispa) -> becomes 112164565 / 92831 (Mnium stellare)
$ws.Range("A9").Value = 112164565
$ws.Range("B9").Value = 92831
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 2362
$ws.Range("F9").Value = "Blek stjärnmossa"
$ws.Range("G9").Value = "Mnium stellare"
$ws.Range("H9").Value = "Hedw."
$ws.Range("Q9").Value = 332935
$ws.Range("R9").Value = 6626957

# Row 10: Taxonsorteringsordning updated
$ws.Range("B10").Value = 95388
